$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'S. GUEYE', 'O. BALSELLS PLAZA', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'S. GUEYE', 'J. ABRIL RAMIRO', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'J. ABRIL RAMIRO', 'S. GUEYE', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'A. GASE SECO', 'K. SKUTYELNIK', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'A. GASE SECO', 'S. GUEYE', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'A. GASE SECO', 'E. LLACER SIMLAT', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'E. LLACER SIMLAT', 'S. GUEYE', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'AUTOCARES RODRÍGUEZ', 'F. ESPINOSA', 'D. VAL GUTIERREZ', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'A. JIMÉNEZ BARRETO', 'O. BALSELLS PLAZA', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'AUTOCARES RODRÍGUEZ', 'D. VAL GUTIERREZ', 'R. DEUS VERA', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'O. BALSELLS PLAZA', 'J. ABRIL RAMIRO', 2)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'AUTOCARES RODRÍGUEZ', 'F. ESPINOSA', 'L. RODRIGUEZ SANCHEZ', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'O. BALSELLS PLAZA', 'J. LARREA IBARBURU', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'O. BALSELLS PLAZA', 'D. VIERA LOPEZ', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'AUTOCARES RODRÍGUEZ', 'R. DEUS VERA', 'D. VAL GUTIERREZ', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'AUTOCARES RODRÍGUEZ', 'D. ROPERO MORA', 'L. RODRIGUEZ SANCHEZ', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'AUTOCARES RODRÍGUEZ', 'P. BARCALA BELTRAN', 'D. VAL GUTIERREZ', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'AUTOCARES RODRÍGUEZ', 'B. BAILLO HERNÁNDEZ', 'D. VAL GUTIERREZ', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'K. SKUTYELNIK', 'S. GUEYE', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'AUTOCARES RODRÍGUEZ', 'L. RODRIGUEZ SANCHEZ', 'D. VAL GUTIERREZ', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'AUTOCARES RODRÍGUEZ', 'D. VAL GUTIERREZ', 'L. RODRIGUEZ SANCHEZ', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'J. ABRIL RAMIRO', 'J. LARREA IBARBURU', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'AUTOCARES RODRÍGUEZ', 'L. RODRIGUEZ SANCHEZ', 'F. ESPINOSA', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'CB ARIDANE', 'E. LLACER SIMLAT', 'A. JIMÉNEZ BARRETO', 1)
    ,@('Liga Regular "B-B"', 1, 'AUTOCARES RODRÍGUEZ vs CB ARIDANE', 'AUTOCARES RODRÍGUEZ', 'S. VAL GUTIERREZ', 'R. DEUS VERA', 1)
)

$startRow = 128
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

Write-Host "Done writing rows"
